# Generate excel file with styles and datas
#
# The sheet originally holds 4 data rows (Empresa 1..4) with no header row.
# This script inserts a new header row on top of the data and styles each
# header cell with a (green) fill, mirroring the column names used
# downstream: nome_fantasia, razao_social, cnpj, endereco, numero, bairro,
# cidade.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data down one row and create a fresh header row at the
# top of the sheet.
$ws.Rows("1:1").Insert()

$headers = @(
    "nome_fantasia",
    "razao_social",
    "cnpj",
    "endereco",
    "numero",
    "bairro",
    "cidade"
)

$green = 65280  # RGB(0, 255, 0)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headers[$i]
    $cell.Interior.Color = $green
}
